# Auto commit update: refresh Metrics values and move active-cell selections.
$wb = $excel.ActiveWorkbook

# --- Metrics sheet: update the underlying raw metric values ---
$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value = 243899.57000000007
$metrics.Range("B3").Value = 214790.13
$metrics.Range("B4").Value = 75336.460000000006
$metrics.Range("B5").Value = 9937
$metrics.Range("B6").Value = 5040145.3200000022
$metrics.Range("B7").Value = 4256866.8100000005
$metrics.Range("B8").Value = 1482296.29
$metrics.Range("B9").Value = 196144
$metrics.Range("B10").Value = 33505526.31000001
$metrics.Range("B11").Value = 31532141.969999999
$metrics.Range("B12").Value = 11764018.33
$metrics.Range("B13").Value = 1293774

# Move the Metrics sheet selection to match the captured state.
$metrics.Range("G17").Select()

# --- today sheet: move the active selection (values there are formulas
#     pulling from Metrics, so they recompute automatically) ---
$today = $wb.Worksheets.Item("today")
$today.Activate()
$today.Range("E7").Select()
